$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 51: coin renamed from EnergySwap to RenderToken, link updated
$ws.Range("B51").Value2 = "RenderToken"
$ws.Range("C51").Value2 = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = "54.574.96"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value2 = "  -6.84%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = "2.430.31"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value2 = "  -10.41%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value2 = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "469.69"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value2 = "  -6.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "133.05"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value2 = "  -5.74%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = "0.994"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value2 = "  -0.42%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "0.495"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value2 = "  -6.62%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "2.448.38"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value2 = "  -10.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "0.0962"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value2 = "  -8.58%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value2 = "  -12.16%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "0.317"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value2 = "  -9.22%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value2 = "  -3.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "2.839.26"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value2 = "  -11.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "54.745.15"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value2 = "  -6.57%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "0.0000135"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value2 = "  -0.64%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "19.84"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value2 = "  -8.66%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "2.432.74"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value2 = "  -10.75%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "4.22"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value2 = "  -11.47%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "313.38"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value2 = "  -8.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "9.55"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value2 = "  -13.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "0.994"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value2 = "  -0.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "5.68"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value2 = "  +0.80%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "5.42"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value2 = "  -13.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "56.67"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value2 = "  +0.95%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "0.388"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value2 = "  -9.22%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "0.158"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value2 = "  -9.49%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "2.526.06"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value2 = "  -11.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "7.19"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value2 = "  -4.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "0.996"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value2 = "  -0.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "0.0₃0721"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value2 = "  -13.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "146.77"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value2 = "  -2.67%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "17.85"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value2 = "  -7.28%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "5.04"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value2 = "  -7.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "3.59"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value2 = "  -14.83%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value2 = "  -6.45%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "0.808"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value2 = "  -14.72%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value2 = "  -0.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "33.06"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value2 = "  -8.24%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "0.598"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value2 = "  -0.26%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "0.0527"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value2 = "  -6.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "3.27"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value2 = "  -8.45%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "1.25"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value2 = "  -10.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "10.10"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value2 = "  -2.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "1.935.40"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value2 = "  -11.70%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value2 = "  -0.20%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "0.0220"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value2 = "  -3.51%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "235.79"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value2 = "  +4.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "4.23"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value2 = "  -10.66%  "
